$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "IST" (actual) hours in column D for the finalized tasks
$ws.Range("D16").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("D19").Value = 1
$ws.Range("D21").Value = 0.5
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 2
$ws.Range("D25").Value = 0.5

# Update the "done" percentage column E to reflect completion
$ws.Range("E16").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("E25").Value = 1

# Update the sheet view: move the selection (also resets the scrolled
# top-left cell back to the default)
$ws.Range("D30").Select()
